# ---------------------------------------------------------------------------
# 1. Merge the two adjacent runs "保存进度 " + "保存当前的工作进度。会分别对暂
#    存区和工作区的状态进行保存。" (identical run formatting) into a single run
#    holding the concatenated text.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute("保存进度 保存当前的工作进度。会分别对暂存区和工作区的状态进行保存。", `
              $true, $false, $false, $false, $false, $true, 1, $false, `
              "保存进度 保存当前的工作进度。会分别对暂存区和工作区的状态进行保存。", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Append four new paragraphs at the end of the document:
#      a) an empty paragraph
#      b) a bold, purple (9B00D3) heading "Git 冲突解决"
#      c) a bold, purple (9B00D3) comment-block line
#      d) a trailing empty paragraph
#
# Insert all four paragraph breaks first -- while the formatting at the
# document's insertion point is still plain (inherited from the preceding
# "快速标志删除..." paragraph: not bold, color=auto) -- and only afterwards
# apply the bold/purple formatting to the two inner (heading) paragraphs.
# That way the leading and trailing empty paragraphs keep the plain format.
# ---------------------------------------------------------------------------
$purple = 13828251   # RGB 9B00D3 expressed as a BGR OLE color value

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

$total = $d.Paragraphs.Count
$pHeading = $d.Paragraphs.Item($total - 2)
$pComment = $d.Paragraphs.Item($total - 1)

$pHeading.Range.InsertBefore("Git 冲突解决")
$pHeading.Range.Font.Bold = $true
$pHeading.Range.Font.Color = $purple

$pComment.Range.InsertBefore("/*********************************************/")
$pComment.Range.Font.Bold = $true
$pComment.Range.Font.Color = $purple
